$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: After "Is the paper well written? How do you know?" add a new run
# with the "For week 2 & later..." sentence, and move the _GoBack bookmark
# here. This also absorbs (removes) the following empty indented paragraph.
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("Is the paper well written? How do you know?", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$insertPos = $r.End

$ins = $d.Range($insertPos, $insertPos)
$xmlNewPara = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="6C7DE16D" w14:textId="77777777" w:rsidR="00D4701F" w:rsidRDefault="00D4701F" w:rsidP="00D4701F"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve"> For week 2 &amp; later, use this space to practice headlines &amp; summaries of the articles via tweets.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
$ins.InsertXML($xmlNewPara)

# Merge the original paragraph (ending at $insertPos) forward into the
# paragraph we just inserted, so the combined paragraph keeps the newly
# inserted paragraph's identity/pPr (ListParagraph/numPr) that we set above
# to match the original.
$pmark1 = $d.Range($insertPos, $insertPos + 1)
$pmark1.Delete()

# Locate the end of the text we just inserted, then remove the very next
# (now-redundant) empty "ind left=360" paragraph by deleting its own
# paragraph mark - this merges its (empty) content into the paragraph
# that follows it, which is visually/structurally identical, net effect:
# one fewer empty placeholder paragraph, matching the diff.
$r2 = $d.Content
$r2.Find.Execute("For week 2 & later, use this space to practice headlines & summaries of the articles via tweets.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$afterNewText = $r2.End
$pmark2 = $d.Range($afterNewText + 1, $afterNewText + 2)
$pmark2.Delete()

# ---------------------------------------------------------------------------
# Change 2: Remove the _GoBack bookmark that used to sit at the end of the
# "Additional Resources:" paragraph (after "... science communication
# pieces."), since it has effectively moved to the location handled above.
# ---------------------------------------------------------------------------
$p3 = $d.Content
$p3.Find.Execute("Additional Resources:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$pstart = $p3.Start

$p4 = $d.Content
$p4.Find.Execute("This is helpful to consider for your science communication pieces.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$pend = $p4.End

$target = $d.Range($pstart, $pend + 1)
$xmlRebuild = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="6EADC32A" w14:textId="50813E8A" w:rsidR="00D4701F" w:rsidRDefault="00D4701F" w:rsidP="00C96AA6"><w:r w:rsidRPr="00C96AA6"><w:rPr><w:b/></w:rPr><w:t>Additional Resources:</w:t></w:r><w:r w:rsidR="00C96AA6"><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>What are the basic concepts that you need to know to understand the science presented in your paper?</w:t></w:r><w:r w:rsidR="00C96AA6"><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>What other information or resources would help you better understand the paper?</w:t></w:r><w:r w:rsidR="00BF2B03"><w:t xml:space="preserve"> This is helpful to consider for your science communication pieces.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="1B2525C2" w14:textId="77777777" w:rsidR="00D4701F" w:rsidRDefault="00D4701F" w:rsidP="00D4701F"/>'
$target.InsertXML($xmlRebuild)

$p5 = $d.Content
$p5.Find.Execute("This is helpful to consider for your science communication pieces.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$afterP5 = $p5.End
$pmark3 = $d.Range($afterP5 + 1, $afterP5 + 2)
$pmark3.Delete()

Write-Output "done"
